# Budgeting Spreadsheet Heather_Andre.xlsx
# Commit message: "Changed the rent amount"
#
# Andre's monthly Rent line (row 9, columns B:M) goes from 1780 to 1800.
# All downstream totals / percentages (rows 13, 14, 42, 43, 45, 46 on the
# Andre sheet, plus the Combined sheet's shared-formula bookkeeping) are
# plain formulas, so they recalculate automatically once this one input
# changes -- no other cell needs to be touched directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Andre")

$ws.Range("B9:M9").Value = 1800

# Leave the selection where the author's last save left it.
$ws.Activate()
$ws.Range("O10").Select() | Out-Null
